# Slide 1 ("Architecture"): the small rectangle that used to read "Sensors"
# is being repurposed as the first modeled Entity box, so relabel it
# "PMEntity" (per commit: "first Entity complete with Data Structure and
# create() function").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 8")
$shp.TextFrame.TextRange.Text = "PMEntity"
